# Update values in the "all_passive_force at subject max ROM" sheet.
# Header row (row 1, columns B:E) - subject counts
# Data rows 2:3 (CON/STR) - force values for the same columns

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 values
$ws.Range("B2").Value = 360.51384934303042
$ws.Range("C2").Value = 813.75693527432168
$ws.Range("D2").Value = 367.55330419549136
$ws.Range("E2").Value = 588.96701264845149

# Row 3 values
$ws.Range("B3").Value = 444.87390981478114
$ws.Range("C3").Value = 436.76752806676581
$ws.Range("D3").Value = 359.55529361567272
$ws.Range("E3").Value = 379.12576171727716

# Update the selected range shown in the sheet view to reflect the edited region
$ws.Range("B1:E3").Select()
